# Overview_routes_dbase.pptx — "Added Maps functionalities and connected
# maps to dbase": shift most shapes down to make room, nudge two shapes
# that moved further (both x and y), and add a new connector arrow.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shapes that shift straight down (Left unchanged) ----------------
$sh = $s.Shapes.Item("Rechthoek 3")
$sh.Top = 77.30622
$sh = $s.Shapes.Item("Rechthoek 4")
$sh.Top = 173.99544
$sh = $s.Shapes.Item("Rechthoek 5")
$sh.Top = 272.6017
$sh = $s.Shapes.Item("Rechthoek 6")
$sh.Top = 371.2079
$sh = $s.Shapes.Item("Rechthoek 7")
$sh.Top = 174.954
$sh = $s.Shapes.Item("Rechthoek 9")
$sh.Top = 474.4132
$sh = $s.Shapes.Item("Tekstvak 1")
$sh.Top = 48.225
$sh = $s.Shapes.Item("Tekstvak 2")
$sh.Top = 48.225
$sh = $s.Shapes.Item("Ovaal 12")
$sh.Top = 84.4406
$sh = $s.Shapes.Item("Ovaal 13")
$sh.Top = 138.12173
$sh = $s.Shapes.Item("Ovaal 14")
$sh.Top = 193.7697
$sh = $s.Shapes.Item("Rechte verbindingslijn met pijl 16")
$sh.Top = 104.55465
$sh = $s.Shapes.Item("Rechte verbindingslijn met pijl 18")
$sh.Top = 107.35072
$sh = $s.Shapes.Item("Rechte verbindingslijn met pijl 20")
$sh.Top = 108.8375
$sh = $s.Shapes.Item("Tekstvak 21")
$sh.Top = 145.8727
$sh = $s.Shapes.Item("Tekstvak 22")
$sh.Top = 245.45221
$sh = $s.Shapes.Item("Tekstvak 23")
$sh.Top = 340.539
$sh = $s.Shapes.Item("Rechthoek 24")
$sh.Top = 472.97324
$sh = $s.Shapes.Item("Tekstvak 25")
$sh.Top = 437.60822
$sh = $s.Shapes.Item("Rechte verbindingslijn met pijl 27")
$sh.Top = 203.223
$sh = $s.Shapes.Item("Rechte verbindingslijn met pijl 29")
$sh.Top = 203.223
$sh = $s.Shapes.Item("Tekstvak 30")
$sh.Top = 252.94142
$sh = $s.Shapes.Item("Rechte verbindingslijn met pijl 32")
$sh.Top = 503.03152
$sh = $s.Shapes.Item("Tekstvak 34")
$sh.Top = 505.44843
$sh = $s.Shapes.Item("Rechte verbindingslijn met pijl 36")
$sh.Top = 217.39261
$sh = $s.Shapes.Item("Rechthoek 37")
$sh.Top = 388.0225

# --- Shapes that moved both horizontally and vertically ---------------
$sh = $s.Shapes.Item("Rechthoek 8")
$sh.Left = 20.1661
$sh.Top = 353.0104

$sh = $s.Shapes.Item("Rechthoek 11")
$sh.Left = 20.0679
$sh.Top = 268.9538

# --- New connector: "Rechte verbindingslijn met pijl 15" ---------------
# Consume the next free shape id (11) with a scratch duplicate so the
# real new connector lands on id 16, matching the authored deck, then
# discard the scratch shape.
$scratch = $s.Shapes.Item("Rechte verbindingslijn met pijl 16").Duplicate()
$scratch.Delete()

$newConn = $s.Shapes.Item("Rechte verbindingslijn met pijl 16").Duplicate()
$newConn.Name = "Rechte verbindingslijn met pijl 15"
$newConn.Left = 177.32033
$newConn.Top = 302.65994
$newConn.Width = 268.6748
$newConn.Height = 0
